$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Raw OOXML widths: C 84->61, D 38->78, F 16->17, H 42->29
# The COM ColumnWidth property differs from the raw stored width by a
# small constant padding offset (~0.83 for this workbook's default font),
# so subtract that offset to land exactly on the target stored width.
$ws.Columns.Item(3).ColumnWidth = 60.17
$ws.Columns.Item(4).ColumnWidth = 77.17
$ws.Columns.Item(6).ColumnWidth = 16.17
$ws.Columns.Item(8).ColumnWidth = 28.17

# --- Row 2 loses its "Yes"/yellow-highlight premium formatting ---
# Revert E2 back to the plain default style (removes the fillId=3 xf).
$ws.Range("E2").Style = "Normal"

# --- Opportunity IDs must stay text, not auto-convert to numbers ---
# (Source data stores these as plain strings, e.g. "1326700".)
$ws.Range("A2:A10").NumberFormat = "@"

# --- Replace the data rows (2-5) and append new rows (6-10) ---
$data = @(
    @("1326700", "https://aiesec.org/opportunity/global-talent/1326700", "Accelerate Romania | Sales Specialist", "Timișoara, Romania", "No", "3 applicants", "9 - 12 Weeks", "re:solved"),
    @("1325702", "https://aiesec.org/opportunity/global-talent/1325702", "Guest Relations Executive and Waitress", "Colombo, Sri Lanka", "No", "8 applicants", "3 - 6 Months", "Indian Kitchen PVT LTD"),
    @("1325344", "https://aiesec.org/opportunity/global-talent/1325344", "Customer Service for finance & accounting (German Speaker)", "Santiago de Querétaro, Qro., Mexico", "No", "24 applicants", "3 - 6 Months", "WMP Mexico Advisors"),
    @("1324549", "https://aiesec.org/opportunity/global-talent/1324549", "Sales Account Manager", "Nasr City, Al Manteqah Al Oula, Nasr City, Cairo Governorate 4450113, Egypt", "No", "18 applicants", "9 - 12 Weeks", "M911 Marketing Emer-Agency"),
    @("1324500", "https://aiesec.org/opportunity/global-talent/1324500", "DT Software Engineer Trainee (EU ONLY)", "Brussels, Belgium", "No", "35 applicants", "6 - 18 Months", "UCB"),
    @("1321823", "https://aiesec.org/opportunity/global-talent/1321823", "Sales Responsible at OnurPlas", "Konya, Türkiye", "No", "43 applicants", "6 - 18 Months", "Onur Plastic"),
    @("1310229", "https://aiesec.org/opportunity/global-talent/1310229", "Guest Relations Officer", "Weligama, Sri Lanka", "No", "64 applicants", "3 - 6 Months", "Weligama Cliff"),
    @("1306542", "https://aiesec.org/opportunity/global-talent/1306542", "Sales Responsible", "Ürgüp, Nevşehir, Türkiye", "No", "103 applicants", "6 - 18 Months", "Pink Lotus Jewellery"),
    @("1303648", "https://aiesec.org/opportunity/global-talent/1303648", "Social Media Manager", "Ahangama, Sri Lanka", "No", "134 applicants", "3 - 6 Months", "Surfing Wombats")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $row = $row + 1
}
